$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "Sum total: x1:0.849|x2:1.000|x3:2.000|x4:1.000|x5:1.000|x6:1.000|x7:1.000|x8:1.000|x9:1.000|x10:1.151|x11:1.000|x12:1.000|x13:1.000|x14:1.000|x15:1.000|x16:1.000|x17:0.000"
$ws.Range("I2").Value = "1.000, 0.347"
$ws.Range("J2").Value = "88.710, 90.323"
$ws.Range("K2").Value = "100.000, 100.000"
